# Refresh the cryptos price table (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (e.g. "71.167.84", "1.00") that must stay
# text, not get auto-coerced to numbers and lose formatting -- force the
# whole data range to Text format before writing the new values.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '71.167.84'
$ws.Range("E2").Value = '  +2.59%  '

$ws.Range("D3").Value = '3.562.50'
$ws.Range("E3").Value = '  +5.15%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").Value = '592.15'
$ws.Range("E5").Value = '  +0.89%  '

$ws.Range("D6").Value = '181.86'
$ws.Range("E6").Value = '  +1.09%  '

$ws.Range("D7").Value = '3.557.19'
$ws.Range("E7").Value = '  +5.28%  '

$ws.Range("E8").Value = '  +1.59%  '

$ws.Range("E9").Value = '  +0.01%  '

$ws.Range("D10").Value = '0.206'
$ws.Range("E10").Value = '  +6.11%  '

$ws.Range("D11").Value = '0.603'
$ws.Range("E11").Value = '  +2.12%  '

$ws.Range("D12").Value = '49.65'
$ws.Range("E12").Value = '  +2.44%  '

$ws.Range("E13").Value = '  +2.16%  '

$ws.Range("D14").Value = '695.63'
$ws.Range("E14").Value = '  +2.69%  '

$ws.Range("D15").Value = '4.119.19'
$ws.Range("E15").Value = '  +4.81%  '

$ws.Range("D16").Value = '8.85'
$ws.Range("E16").Value = '  +2.84%  '

$ws.Range("D17").Value = '71.258.35'
$ws.Range("E17").Value = '  +2.68%  '

$ws.Range("D18").Value = '3.536.80'
$ws.Range("E18").Value = '  +4.37%  '

$ws.Range("E19").Value = '  +1.19%  '

$ws.Range("D20").Value = '18.24'
$ws.Range("E20").Value = '  +3.48%  '

$ws.Range("D21").Value = '11.63'
$ws.Range("E21").Value = '  +3.20%  '

$ws.Range("D22").Value = '0.923'
$ws.Range("E22").Value = '  +2.25%  '

$ws.Range("D23").Value = '5.51'
$ws.Range("E23").Value = '  +1.76%  '

$ws.Range("D24").Value = '17.46'
$ws.Range("E24").Value = '  +1.53%  '

$ws.Range("D25").Value = '103.10'
$ws.Range("E25").Value = '  -0.20%  '

$ws.Range("E26").Value = '  +1.55%  '

$ws.Range("D27").Value = '2.77'
$ws.Range("E27").Value = '  +1.65%  '

$ws.Range("D28").Value = '9.88'
$ws.Range("E28").Value = '  +2.40%  '

$ws.Range("D29").Value = '34.62'
$ws.Range("E29").Value = '  +1.62%  '

$ws.Range("D30").Value = '9.01'
$ws.Range("E30").Value = '  +3.60%  '

$ws.Range("E31").Value = '  +4.18%  '

$ws.Range("E32").Value = '  +11.55%  '

$ws.Range("D33").Value = '583.11'
$ws.Range("E33").Value = '  +4.86%  '

$ws.Range("D34").Value = '11.25'
$ws.Range("E34").Value = '  +1.20%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.106'
$ws.Range("E35").Value = '  -0.37%  '

$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D36").Value = '59.33'
$ws.Range("E36").Value = '  +1.79%  '

$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.01%  '

$ws.Range("D38").Value = '3.657.60'
$ws.Range("E38").Value = '  -0.61%  '

$ws.Range("E39").Value = '  +2.83%  '

$ws.Range("D40").Value = '35.85'
$ws.Range("E40").Value = '  +1.60%  '

$ws.Range("B41").Value = 'PEPE'
$ws.Range("C41").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D41").Value = '0.0₃0757'
$ws.Range("E41").Value = '  +8.82%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '3.43'
$ws.Range("E42").Value = '  +5.03%  '

$ws.Range("E43").Value = '  +3.25%  '

$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0435'
$ws.Range("E44").Value = '  +2.85%  '

$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").Value = '0.345'
$ws.Range("E45").Value = '  +1.79%  '

$ws.Range("D46").Value = '3.39'
$ws.Range("E46").Value = '  +2.74%  '

$ws.Range("D47").Value = '2.74'
$ws.Range("E47").Value = '  +2.70%  '

$ws.Range("D48").Value = '1.48'
$ws.Range("E48").Value = '  +4.63%  '

$ws.Range("E49").Value = '  +1.12%  '

$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  -0.05%  '

$ws.Range("D51").Value = '133.82'
$ws.Range("E51").Value = '  +0.82%  '
